$d = $word.ActiveDocument

# Change 1: shorten the "Discussion" section summary bullet.
$d.Content.Find.Execute(
    "The first point in the " + [char]0x201C + "Discussion" + [char]0x201D + " section focused on the comparison between deep learning methods and probabilistic methods. In particular, this point explained the advantages of RNN over HMM or MC.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The first point in the " + [char]0x201C + "Discussion" + [char]0x201D + " section explained the advantages of RNN over probabilistic models.",
    2
)

# Change 2: reword the "Results" section bullet about the PPC abbreviation.
$d.Content.Find.Execute(
    "Modified the " + [char]0x201C + "Results" + [char]0x201D + " section to use the full name of PPC code.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Modified the " + [char]0x201C + "Results" + [char]0x201D + " section to resolve the problem of the abbreviation " + [char]0x201C + "PPC" + [char]0x201D + ".",
    2
)

# Change 3: fix "CT" -> "CHT" (both occurrences) in the Reviewer #4 paragraph.
$d.Content.Find.Execute(
    "examples of CT and CML specific to black adolescents",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "examples of CHT and CML specific to black adolescents",
    2
)

$d.Content.Find.Execute(
    "examples of CT and CML for the pertinent target behaviors",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "examples of CHT and CML for the pertinent target behaviors",
    2
)

# Change 3 (cont.): "However, this study" -> "This study".
$d.Content.Find.Execute(
    "the sequence of words. However, this study represents",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "the sequence of words. This study represents",
    2
)
